# Update "Ngày kết thúc thực tế" (actual end date) for the usecase-related
# rows on Sheet1.
#   Row 9  (Thiết kế mô hình cơ sở dữ liệu)      -> 21/09/2020
#   Row 10 (Thu thập thông tin và nhập vào CSDL) -> 23/09/2020
#   Row 11 (Tìm hiểu ASP.NET Core, CSS, HTML...) -> 21/09/2020

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H9").Value  = 44095
$ws.Range("H10").Value = 44097
$ws.Range("H11").Value = 44095

# Reflect the reviewer's resulting selection state after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H20").Select()
